{"js": "// Add \", Render, Railway\" to the end of the \"Deployment and Hosting\" skills\n// line, and \", Docker\" to the end of the \"Tools\" skills line \u2014 each as its\n// own new run (matching how Word appends freshly-typed text as a distinct\n// run rather than merging into the previous run).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet hostingParagraph = null;\nlet toolsParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  if (hostingParagraph === null && text.indexOf(\"Netlify, Digital Ocean\") !== -1) {\n    hostingParagraph = paragraph;\n  }\n  if (toolsParagraph === null && /Postman,\\s*Figma\\s*$/.test(text)) {\n    toolsParagraph = paragraph;\n  }\n}\n\nif (hostingParagraph) {\n  hostingParagraph.insertText(\", Render, Railway\", Word.InsertLocation.end);\n}\nif (toolsParagraph) {\n  toolsParagraph.insertText(\", Docker\", Word.InsertLocation.end);\n}\n\nawait context.sync();\n", "ps1": "# Add \", Render, Railway\" to the end of the \"Deployment and Hosting\" skills\n# line, and \", Docker\" to the end of the \"Tools\" skills line \u2014 each lands as\n# its own new run, matching how Word appends freshly-typed text as a\n# separate run rather than merging into the previous run's text.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Netlify, Digital Ocean*\") {\n        $p.Range.InsertAfter(\", Render, Railway\")\n    }\n    elseif ($t -like \"*Postman, Figma*\") {\n        $p.Range.InsertAfter(\", Docker\")\n    }\n}\n"}
